$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the first sheet
$ws.Name = "Export as TSV"

# 2. Freeze the header row (row 1) on the main sheet
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

# 3. Add error titles/messages to the existing data validations
$rules = @(
    @{ Range = "I2:I1048576"; Title = "Value must come from list"; Message = "Value must be one of: imaging." },
    @{ Range = "J2:J1048576"; Title = "Value must come from list"; Message = "Value must be one of: MxIF." },
    @{ Range = "K2:K1048576"; Title = "Value must come from list"; Message = "Value must be one of: protein." },
    @{ Range = "L2:L1048576"; Title = "Not a boolean"; Message = 'The values in this column must be "TRUE" or "FALSE".' },
    @{ Range = "O2:O1048576"; Title = "Not a number"; Message = "The values in this column must be numbers." },
    @{ Range = "P2:P1048576"; Title = "Value must come from list"; Message = "Value must be one of: nm / um." },
    @{ Range = "Q2:Q1048576"; Title = "Not a number"; Message = "The values in this column must be numbers." },
    @{ Range = "R2:R1048576"; Title = "Value must come from list"; Message = "Value must be one of: nm / um." },
    @{ Range = "S2:S1048576"; Title = "Not an integer"; Message = "The values in this column must be integers." },
    @{ Range = "T2:T1048576"; Title = "Not an integer"; Message = "The values in this column must be integers." }
)

foreach ($rule in $rules) {
    $validation = $ws.Range($rule.Range).Validation
    $validation.ErrorTitle = $rule.Title
    $validation.ErrorMessage = $rule.Message
}
